{"js": "// Replace the date line and every three-digit x one-digit multiplication\n// equation in the answer table with the new values from the target revision.\n// Every <w:t> run in the document changes (1 title run + 25 table-cell runs),\n// so we do an exact-text search/replace for each old -> new pair. All the\n// \"old\" strings are unique substrings across the document, so matching on\n// them unambiguously targets the correct run.\nconst replacements = [\n  [\"2024-08-23 Friday\", \"2024-08-24 Saturday\"],\n  [\"601\u00d72=1202\", \"797\u00d77=5579\"],\n  [\"381\u00d78=3048\", \"822\u00d76=4932\"],\n  [\"464\u00d73=1392\", \"654\u00d75=3270\"],\n  [\"447\u00d78=3576\", \"285\u00d79=2565\"],\n  [\"296\u00d77=2072\", \"544\u00d72=1088\"],\n  [\"825\u00d74=3300\", \"312\u00d72=624\"],\n  [\"407\u00d73=1221\", \"349\u00d78=2792\"],\n  [\"452\u00d72=904\", \"710\u00d72=1420\"],\n  [\"294\u00d75=1470\", \"379\u00d73=1137\"],\n  [\"199\u00d76=1194\", \"683\u00d77=4781\"],\n  [\"589\u00d75=2945\", \"148\u00d77=1036\"],\n  [\"122\u00d75=610\", \"646\u00d76=3876\"],\n  [\"162\u00d78=1296\", \"976\u00d79=8784\"],\n  [\"299\u00d76=1794\", \"787\u00d74=3148\"],\n  [\"781\u00d76=4686\", \"708\u00d74=2832\"],\n  [\"962\u00d75=4810\", \"142\u00d79=1278\"],\n  [\"952\u00d73=2856\", \"246\u00d76=1476\"],\n  [\"795\u00d72=1590\", \"247\u00d77=1729\"],\n  [\"970\u00d77=6790\", \"283\u00d73=849\"],\n  [\"974\u00d78=7792\", \"494\u00d79=4446\"],\n  [\"669\u00d76=4014\", \"678\u00d74=2712\"],\n  [\"588\u00d75=2940\", \"241\u00d73=723\"],\n  [\"595\u00d75=2975\", \"591\u00d74=2364\"],\n  [\"109\u00d77=763\", \"681\u00d78=5448\"],\n  [\"864\u00d77=6048\", \"696\u00d76=4176\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every three-digit x one-digit multiplication\n# equation in the answer table with the new values from the target revision.\n# Every run of text in the document changes (1 title run + 25 table-cell\n# runs), so we do an exact-text Find/Replace for each old -> new pair over\n# the whole document range. All the \"old\" strings are unique substrings\n# across the document, so matching on them unambiguously targets the\n# correct run, and MatchCase keeps the match exact.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-08-23 Friday\", \"2024-08-24 Saturday\"),\n    @(\"601\u00d72=1202\", \"797\u00d77=5579\"),\n    @(\"381\u00d78=3048\", \"822\u00d76=4932\"),\n    @(\"464\u00d73=1392\", \"654\u00d75=3270\"),\n    @(\"447\u00d78=3576\", \"285\u00d79=2565\"),\n    @(\"296\u00d77=2072\", \"544\u00d72=1088\"),\n    @(\"825\u00d74=3300\", \"312\u00d72=624\"),\n    @(\"407\u00d73=1221\", \"349\u00d78=2792\"),\n    @(\"452\u00d72=904\", \"710\u00d72=1420\"),\n    @(\"294\u00d75=1470\", \"379\u00d73=1137\"),\n    @(\"199\u00d76=1194\", \"683\u00d77=4781\"),\n    @(\"589\u00d75=2945\", \"148\u00d77=1036\"),\n    @(\"122\u00d75=610\", \"646\u00d76=3876\"),\n    @(\"162\u00d78=1296\", \"976\u00d79=8784\"),\n    @(\"299\u00d76=1794\", \"787\u00d74=3148\"),\n    @(\"781\u00d76=4686\", \"708\u00d74=2832\"),\n    @(\"962\u00d75=4810\", \"142\u00d79=1278\"),\n    @(\"952\u00d73=2856\", \"246\u00d76=1476\"),\n    @(\"795\u00d72=1590\", \"247\u00d77=1729\"),\n    @(\"970\u00d77=6790\", \"283\u00d73=849\"),\n    @(\"974\u00d78=7792\", \"494\u00d79=4446\"),\n    @(\"669\u00d76=4014\", \"678\u00d74=2712\"),\n    @(\"588\u00d75=2940\", \"241\u00d73=723\"),\n    @(\"595\u00d75=2975\", \"591\u00d74=2364\"),\n    @(\"109\u00d77=763\", \"681\u00d78=5448\"),\n    @(\"864\u00d77=6048\", \"696\u00d76=4176\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
